$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the D2:E51 range to Text format first so numeric-looking strings
# (e.g. "1.00", "132.67") are stored verbatim instead of being normalized
# to numbers by the automatic type inference on Range.Value assignment.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '64.691.78'
$ws.Range("E2").Value = '  -1.81%  '
$ws.Range("D3").Value = '3.514.57'
$ws.Range("E3").Value = '  -2.55%  '
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = '586.75'
$ws.Range("E5").Value = '  -3.08%  '
$ws.Range("D6").Value = '132.67'
$ws.Range("E6").Value = '  -3.13%  '
$ws.Range("D7").Value = '3.514.56'
$ws.Range("E7").Value = '  -2.54%  '
$ws.Range("D9").Value = '0.492'
$ws.Range("E9").Value = '  -1.30%  '
$ws.Range("E10").Value = '  -0.81%  '
$ws.Range("D11").Value = '7.33'
$ws.Range("E11").Value = '  +1.31%  '
$ws.Range("E12").Value = '  -1.22%  '
$ws.Range("D13").Value = '4.112.93'
$ws.Range("E13").Value = '  -1.98%  '
$ws.Range("D14").Value = '27.79'
$ws.Range("E14").Value = '  -1.19%  '
$ws.Range("E15").Value = '  -3.71%  '
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("D17").Value = '3.512.63'
$ws.Range("E17").Value = '  -2.11%  '
$ws.Range("D18").Value = '64.667.48'
$ws.Range("E18").Value = '  -0.48%  '
$ws.Range("E19").Value = '  -1.32%  '
$ws.Range("D20").Value = '14.21'
$ws.Range("E20").Value = '  -3.13%  '
$ws.Range("E21").Value = '  -4.36%  '
$ws.Range("D22").Value = '391.67'
$ws.Range("E22").Value = '  -1.61%  '
$ws.Range("E23").Value = '  -1.77%  '
$ws.Range("D24").Value = '3.655.34'
$ws.Range("E24").Value = '  -2.44%  '
$ws.Range("D25").Value = '73.74'
$ws.Range("E25").Value = '  -1.26%  '
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("E27").Value = '  -5.60%  '
$ws.Range("E28").Value = '  -5.75%  '
$ws.Range("D29").Value = '7.49'
$ws.Range("E29").Value = '  -8.26%  '
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("E31").Value = '  -5.40%  '
$ws.Range("D32").Value = '8.23'
$ws.Range("E32").Value = '  -5.81%  '
$ws.Range("D33").Value = '3.515.07'
$ws.Range("E33").Value = '  -2.38%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").Value = '24.05'
$ws.Range("E35").Value = '  -2.47%  '
$ws.Range("E36").Value = '  -0.91%  '
$ws.Range("D37").Value = '5.29'
$ws.Range("E37").Value = '  -0.73%  '
$ws.Range("E38").Value = '  +0.17%  '
$ws.Range("D39").Value = '170.99'
$ws.Range("E39").Value = '  -0.30%  '
$ws.Range("D40").Value = '6.99'
$ws.Range("E40").Value = '  -1.66%  '
$ws.Range("D41").Value = '0.0807'
$ws.Range("E41").Value = '  -3.29%  '
$ws.Range("D42").Value = '26.56'
$ws.Range("E42").Value = '  +0.27%  '
$ws.Range("E43").Value = '  -3.52%  '
$ws.Range("E44").Value = '  +0.46%  '
$ws.Range("D45").Value = '42.13'
$ws.Range("E45").Value = '  -2.87%  '
$ws.Range("D46").Value = '1.21'
$ws.Range("E46").Value = '  -2.28%  '
$ws.Range("E47").Value = '  -2.98%  '
$ws.Range("E48").Value = '  -3.10%  '
$ws.Range("D49").Value = '2.448.94'
$ws.Range("E49").Value = '  -0.94%  '
$ws.Range("D50").Value = '6.90'
$ws.Range("E50").Value = '  -2.65%  '
$ws.Range("D51").Value = '0.905'
$ws.Range("E51").Value = '  +2.51%  '

# Restore the default (un-styled) cell style so no stray formatting is
# left behind by the temporary Text number format applied above.
$dataRange.Style = "Normal"
